# Update calibration data with new costs for rows 100-107 and 114-115
# in the "strategy_id-0" sheet. Columns J:AS hold a uniform value per row
# (one value repeated across all year columns); set each row's J:AS range
# to its new uniform value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    100 = 0
    101 = 54688.1154
    102 = 512182.4956
    103 = 87409.98522
    104 = 0
    105 = 38443.70824
    106 = 56340.07452
    107 = 257522.7288
    114 = 40.88976632
    115 = 900703.6347000001
}

foreach ($row in $updates.Keys) {
    $value = $updates[$row]
    $range = $ws.Range("J$row`:AS$row")
    $range.Value = $value
}
